$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Fitness column (C) values for rows 2 through 12 to 4357
$ws.Range("C2:C12").Value = 4357
